# textures, spirit consumable, mc wip, consumable panel
# Adds a new "Spirit" consumable column (F) to the Character sheet, mirroring
# the existing Hp/Dmg column (D/E) formatting, and selects the new cell F2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mirror column E's formatting (fill/border/alignment/number format) into the
# new column F for the header + the two data rows that actually carry values.
$ws.Range("E1:E3").Copy()
$ws.Range("F1:F3").PasteSpecial(-4122)

# Header cell F1 repeats the "Hp/\nDmg" header text used by D1/E1.
$ws.Cells.Item(1, 6).Value = $ws.Cells.Item(1, 5).Value2

# F2 holds the new "Spirit" consumable formula text.
$ws.Cells.Item(2, 6).Value = "5 - floor *10"

# F3 repeats the "floor * 2" text used by D3/E3.
$ws.Cells.Item(3, 6).Value = $ws.Cells.Item(3, 4).Value2

# Give column F its own width (23 characters).
$ws.Range("F1").EntireColumn.ColumnWidth = 22.166666666666668

# Move the active selection to the newly-populated cell.
$ws.Range("F2").Select() | Out-Null
